# Add team record (Wins/Losses/Ties) columns to the MIL_2010 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of an existing header cell (bold, bordered, centered)
# onto the three new header cells so they match the rest of the header row.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record values for every data row (rows 2 through 45).
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 30).Value = 77  # AD -> Wins
    $ws.Cells.Item($row, 31).Value = 85  # AE -> Losses
    $ws.Cells.Item($row, 32).Value = 0   # AF -> Ties
}

Write-Output "done"
